# Applies the "Atualizado por script em 12-11-2023 14:45" update:
#   1) Four pairs of adjacent match rows get their match-detail columns
#      (F:V) swapped with each other (the Indice/pais/torneio/temporada/
#      data_partida columns A:E stay put).
#   2) Three brand-new match rows are appended at the bottom of the sheet
#      (rows 106-108, Indice 105-107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the four row pairs whose match details were re-ordered.
$swapPairs = @(
    @(21, 22),
    @(50, 51),
    @(74, 75),
    @(89, 90)
)

foreach ($pair in $swapPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("F$rowA`:V$rowA")
    $rangeB = $ws.Range("F$rowB`:V$rowB")

    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()

    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

# 2) Append the three new match rows (106-108) at the end of the sheet.
$templateRow = 105

$newRows = @(
    @{ Row=106; A=105; E=45242.47916666666; F="Boluspor";    G=2; H="Sakaryaspor";   I=3;
       J=2.79; K="05/11/2023 14:12"; L=2.69; M="12/11/2023 11:01";
       N=3.15; O="05/11/2023 14:12"; P=3.01; Q="12/11/2023 10:41";
       R=2.65; S="05/11/2023 14:12"; T=2.95; U="12/11/2023 11:01";
       V="https://www.betexplorer.com/football/turkey/1-lig/boluspor-sakaryaspor/IcM3n7Ip/" },
    @{ Row=107; A=106; E=45242.47916666666; F="Manisa FK";   G=0; H="Bandirmaspor";  I=2;
       J=2.18; K="05/11/2023 11:42"; L=2.54; M="12/11/2023 11:21";
       N=3.41; O="05/11/2023 11:42"; P=3.35; Q="12/11/2023 11:25";
       R=3.31; S="05/11/2023 11:42"; T=2.84; U="12/11/2023 11:25";
       V="https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-bandirmaspor/Eq666S9G/" },
    @{ Row=108; A=107; E=45242.58333333334; F="Erzurumspor"; G=3; H="Goztepe";       I=2;
       J=4.1;  K="05/11/2023 14:12"; L=3.85; M="12/11/2023 13:52";
       N=3.37; O="05/11/2023 14:12"; P=3.05; Q="12/11/2023 13:54";
       R=1.93; S="05/11/2023 14:12"; T=2.19; U="12/11/2023 13:54";
       V="https://www.betexplorer.com/football/turkey/1-lig/erzurumspor-fk-goztepe/tW7b8lu4/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Clone formatting (bold/border style on A, date format on E) from the
    # last existing data row so the new rows match the sheet's styling.
    $ws.Range("A$templateRow").Copy($ws.Range("A$r"))
    $ws.Range("E$templateRow").Copy($ws.Range("E$r"))

    $ws.Range("A$r").Value = $nr.A
    $ws.Range("B$r").Value = "turkey"
    $ws.Range("C$r").Value = "1-lig"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $nr.E
    $ws.Range("F$r").Value = $nr.F
    $ws.Range("G$r").Value = $nr.G
    $ws.Range("H$r").Value = $nr.H
    $ws.Range("I$r").Value = $nr.I
    $ws.Range("J$r").Value = $nr.J
    $ws.Range("K$r").Value = $nr.K
    $ws.Range("L$r").Value = $nr.L
    $ws.Range("M$r").Value = $nr.M
    $ws.Range("N$r").Value = $nr.N
    $ws.Range("O$r").Value = $nr.O
    $ws.Range("P$r").Value = $nr.P
    $ws.Range("Q$r").Value = $nr.Q
    $ws.Range("R$r").Value = $nr.R
    $ws.Range("S$r").Value = $nr.S
    $ws.Range("T$r").Value = $nr.T
    $ws.Range("U$r").Value = $nr.U
    $ws.Range("V$r").Value = $nr.V
}

Write-Output "Edit complete"
